$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row (row 1 was previously unused/empty) with column titles
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Категория"
$ws.Range("C1").Value = "Пол"
$ws.Range("D1").Value = "Дата рождения"
$ws.Range("E1").Value = "Идентификационная метка"
$ws.Range("F1").Value = "Номер электронного чипа"
$ws.Range("G1").Value = "Кличка"
$ws.Range("H1").Value = "Фото"
$ws.Range("I1").Value = "Владелец"

# Fix row 2 data - correct category/gender for record "Вася"
$ws.Range("B2").Value = "Собака"
$ws.Range("C2").Value = "Мужской"

# Replace row 8 (old record 2002 "Жук") with the data for record 2007 "Зая"
$ws.Range("A8").Value = 2007
$ws.Range("B8").Value = "Собака"
$ws.Range("C8").Value = "Женский"
# "01.12.2021" is ambiguous with a real date (day<=12), so route it through a
# text formula and paste-as-values to keep it a literal text cell.
$ws.Range("D8").Formula = '="01.12.2021"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "Зая"
$ws.Range("I8").Value = "Ромашка"

# Replace row 9 (old record 2007 "Зая") with the data for record 3002 "Жук"
$ws.Range("A9").Value = 3002
$ws.Range("B9").Value = "Собака"
$ws.Range("C9").Value = "Мужской"
$ws.Range("D9").Value = "16.12.2018"
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = "Жук"
$ws.Range("I9").Value = "Sapov EV"
